$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: capture current (pre-edit) row data for columns B:AC for all affected rows
$rowData = @{}
$rowData[21] = $ws.Range("B21:AC21").Value2
$rowData[22] = $ws.Range("B22:AC22").Value2
$rowData[27] = $ws.Range("B27:AC27").Value2
$rowData[28] = $ws.Range("B28:AC28").Value2
$rowData[37] = $ws.Range("B37:AC37").Value2
$rowData[38] = $ws.Range("B38:AC38").Value2
$rowData[42] = $ws.Range("B42:AC42").Value2
$rowData[43] = $ws.Range("B43:AC43").Value2
$rowData[53] = $ws.Range("B53:AC53").Value2
$rowData[54] = $ws.Range("B54:AC54").Value2
$rowData[69] = $ws.Range("B69:AC69").Value2
$rowData[70] = $ws.Range("B70:AC70").Value2
$rowData[84] = $ws.Range("B84:AC84").Value2
$rowData[85] = $ws.Range("B85:AC85").Value2
$rowData[86] = $ws.Range("B86:AC86").Value2
$rowData[87] = $ws.Range("B87:AC87").Value2
$rowData[131] = $ws.Range("B131:AC131").Value2
$rowData[132] = $ws.Range("B132:AC132").Value2
$rowData[133] = $ws.Range("B133:AC133").Value2
$rowData[134] = $ws.Range("B134:AC134").Value2
$rowData[135] = $ws.Range("B135:AC135").Value2
$rowData[175] = $ws.Range("B175:AC175").Value2
$rowData[176] = $ws.Range("B176:AC176").Value2
$rowData[223] = $ws.Range("B223:AC223").Value2
$rowData[224] = $ws.Range("B224:AC224").Value2
$rowData[225] = $ws.Range("B225:AC225").Value2
$rowData[226] = $ws.Range("B226:AC226").Value2
$rowData[227] = $ws.Range("B227:AC227").Value2
$rowData[228] = $ws.Range("B228:AC228").Value2
$rowData[229] = $ws.Range("B229:AC229").Value2
$rowData[230] = $ws.Range("B230:AC230").Value2
$rowData[231] = $ws.Range("B231:AC231").Value2
$rowData[232] = $ws.Range("B232:AC232").Value2
$rowData[237] = $ws.Range("B237:AC237").Value2
$rowData[238] = $ws.Range("B238:AC238").Value2
$rowData[239] = $ws.Range("B239:AC239").Value2
$rowData[240] = $ws.Range("B240:AC240").Value2
$rowData[247] = $ws.Range("B247:AC247").Value2
$rowData[248] = $ws.Range("B248:AC248").Value2

# Step 2: write back according to the row permutation mapping
$ws.Range("B21:AC21").Value2 = $rowData[22]
$ws.Range("B22:AC22").Value2 = $rowData[21]
$ws.Range("B27:AC27").Value2 = $rowData[28]
$ws.Range("B28:AC28").Value2 = $rowData[27]
$ws.Range("B37:AC37").Value2 = $rowData[38]
$ws.Range("B38:AC38").Value2 = $rowData[37]
$ws.Range("B42:AC42").Value2 = $rowData[43]
$ws.Range("B43:AC43").Value2 = $rowData[42]
$ws.Range("B53:AC53").Value2 = $rowData[54]
$ws.Range("B54:AC54").Value2 = $rowData[53]
$ws.Range("B69:AC69").Value2 = $rowData[70]
$ws.Range("B70:AC70").Value2 = $rowData[69]
$ws.Range("B84:AC84").Value2 = $rowData[87]
$ws.Range("B85:AC85").Value2 = $rowData[86]
$ws.Range("B86:AC86").Value2 = $rowData[84]
$ws.Range("B87:AC87").Value2 = $rowData[85]
$ws.Range("B131:AC131").Value2 = $rowData[132]
$ws.Range("B132:AC132").Value2 = $rowData[133]
$ws.Range("B133:AC133").Value2 = $rowData[131]
$ws.Range("B134:AC134").Value2 = $rowData[135]
$ws.Range("B135:AC135").Value2 = $rowData[134]
$ws.Range("B175:AC175").Value2 = $rowData[176]
$ws.Range("B176:AC176").Value2 = $rowData[175]
$ws.Range("B223:AC223").Value2 = $rowData[224]
$ws.Range("B224:AC224").Value2 = $rowData[223]
$ws.Range("B225:AC225").Value2 = $rowData[226]
$ws.Range("B226:AC226").Value2 = $rowData[225]
$ws.Range("B227:AC227").Value2 = $rowData[228]
$ws.Range("B228:AC228").Value2 = $rowData[227]
$ws.Range("B229:AC229").Value2 = $rowData[231]
$ws.Range("B230:AC230").Value2 = $rowData[232]
$ws.Range("B231:AC231").Value2 = $rowData[230]
$ws.Range("B232:AC232").Value2 = $rowData[229]
$ws.Range("B237:AC237").Value2 = $rowData[238]
$ws.Range("B238:AC238").Value2 = $rowData[239]
$ws.Range("B239:AC239").Value2 = $rowData[240]
$ws.Range("B240:AC240").Value2 = $rowData[237]
$ws.Range("B247:AC247").Value2 = $rowData[248]
$ws.Range("B248:AC248").Value2 = $rowData[247]

Write-Host "Row permutation edits applied."
